$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. Values are forced to Text format
# ("@" number format) before assignment so Excel does not auto-convert
# numeric-looking strings (prices, percentages) into floating point numbers,
# preserving the exact text representation from the source data feed.
$updates = [ordered]@{
    'D2' = '288.36'
    'E2' = '1.26%'
    'D3' = '29.22'
    'E3' = '1.61%'
    'D4' = '5.075'
    'E4' = '3.39%'
    'D5' = '0.06675'
    'E5' = '2.90%'
    'D6' = '7.339'
    'E6' = '1.78%'
    'B7' = 'FTXToken'
    'C7' = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
    'D7' = '1.348'
    'E7' = '0.72%'
    'B8' = 'MXToken'
    'C8' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
    'D8' = '0.9186'
    'E8' = '0.51%'
    'B9' = 'WazirX'
    'C9' = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
    'D9' = '0.1587'
    'E9' = '3.02%'
    'B10' = 'LiechtensteinCryptoassetsExchange'
    'C10' = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
    'D10' = '0.06804'
    'E10' = '6.15%'
    'B11' = 'MandalaExchangeToken'
    'C11' = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
    'D11' = '0.07668'
    'E11' = '1.39%'
    'B12' = 'BitrueCoin'
    'C12' = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
    'D12' = '0.02936'
    'E12' = '-1.49%'
    'B13' = 'BitMartToken'
    'C13' = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
    'D13' = '0.08989'
    'E13' = '0.20%'
    'B14' = 'BitForexToken'
    'C14' = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
    'D14' = '0.001565'
    'E14' = '-2.53%'
    'B15' = 'CoinExToken'
    'C15' = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
    'D15' = '0.04513'
    'E15' = '0.84%'
    'B16' = 'One'
    'C16' = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
    'D16' = '0.0006467'
    'E16' = '-1.44%'
    'B17' = 'TigerCash'
    'C17' = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
    'D17' = '0.006253'
    'E17' = '3.35%'
    'B18' = 'LEO'
    'C18' = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
    'D18' = '3.442'
    'E18' = '-0.51%'
    'B19' = 'GateToken'
    'C19' = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
    'D19' = '3.405'
    'E19' = '0.87%'
    'D20' = '2.218'
    'E20' = '-1.06%'
    'E21' = '2.05%'
    'E22' = '-2.47%'
    'D23' = '4.065'
    'E23' = '1.33%'
    'E24' = '1.69%'
    'D25' = '0.001192'
    'E25' = '0.28%'
    'D26' = '0.004112'
    'E26' = '-4.87%'
    'E27' = '1.45%'
    'D40' = '0.04221'
    'E40' = '1.74%'
    'D41' = '0.006725'
    'E41' = '0.68%'
    'E42' = '0.60%'
    'E43' = '-3.84%'
    'D44' = '0.01338'
    'E44' = '13.50%'
    'D45' = '0.00005712'
    'E45' = '6.29%'
    'D46' = '1.974'
    'E46' = '8.52%'
    'E47' = '-29.42%'
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
